# Anonymize "fedcore" -> "approach" in the header rows and give the
# C1/D1 (and F1/G1 on the computational sheet) header cells their own
# thin-border styling, matching the author's re-export of the results.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Build the two new cell styles exactly once (on sheet1), then
# --- replicate them with Copy/PasteSpecial so the engine reuses the
# --- same style record instead of re-deriving it (and leaving stray
# --- intermediate style entries) for every target cell.

# Style A: top + bottom thin border only (no left/right) -> used by C1
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = 1
$c1.Borders.Item(7).LineStyle = -4142
$c1.Borders.Item(10).LineStyle = -4142

# Style B: top + bottom + right thin border (no left) -> used by D1
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.LineStyle = 1
$d1.Borders.Item(7).LineStyle = -4142

# Apply the same two styles to the matching header cells on sheet2
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty inline-string cell G5 on the computational sheet
$ws2.Range("G5").ClearContents()
